$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original Text storage so numeric-looking
# strings (e.g. "0.999", "68.116.67") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '68.116.67'
$ws.Range("E2").Value = '  -5.09%  '

$ws.Range("D3").Value = '3.704.78'
$ws.Range("E3").Value = '  -4.72%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '584.73'
$ws.Range("E5").Value = '  -2.10%  '

$ws.Range("D6").Value = '180.29'
$ws.Range("E6").Value = '  +7.75%  '

$ws.Range("D7").Value = '3.696.45'
$ws.Range("E7").Value = '  -5.00%  '

$ws.Range("D8").Value = '0.629'
$ws.Range("E8").Value = '  -5.92%  '

$ws.Range("D9").Value = '0.998'
$ws.Range("E9").Value = '  -0.17%  '

$ws.Range("D10").Value = '0.712'
$ws.Range("E10").Value = '  -6.42%  '

$ws.Range("D11").Value = '0.163'
$ws.Range("E11").Value = '  -8.74%  '

$ws.Range("D12").Value = '54.03'
$ws.Range("E12").Value = '  +0.26%  '

$ws.Range("D13").Value = '0.0000291'
$ws.Range("E13").Value = '  -9.62%  '

$ws.Range("D14").Value = '10.42'
$ws.Range("E14").Value = '  -8.09%  '

$ws.Range("D15").Value = '4.353.43'
$ws.Range("E15").Value = '  -3.67%  '

$ws.Range("D16").Value = '3.689.04'
$ws.Range("E16").Value = '  -5.44%  '

$ws.Range("D17").Value = '19.46'
$ws.Range("E17").Value = '  -7.19%  '

$ws.Range("E18").Value = '  -2.70%  '

$ws.Range("D19").Value = '12.83'
$ws.Range("E19").Value = '  -7.75%  '

$ws.Range("E20").Value = '  -7.40%  '

$ws.Range("D21").Value = '67.695.09'
$ws.Range("E21").Value = '  -5.60%  '

$ws.Range("D22").Value = '408.37'
$ws.Range("E22").Value = '  -6.08%  '

$ws.Range("E23").Value = '  -5.10%  '

$ws.Range("D24").Value = '88.46'
$ws.Range("E24").Value = '  -6.10%  '

$ws.Range("E25").Value = '  -8.01%  '

$ws.Range("D26").Value = '12.77'
$ws.Range("E26").Value = '  -7.71%  '

$ws.Range("D27").Value = '10.98'
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").Value = '3.87'
$ws.Range("E28").Value = '  -7.06%  '

$ws.Range("D29").Value = '6.06'
$ws.Range("E29").Value = '  +2.22%  '

$ws.Range("D30").Value = '9.51'
$ws.Range("E30").Value = '  -6.61%  '

$ws.Range("D31").Value = '32.52'
$ws.Range("E31").Value = '  -7.39%  '

$ws.Range("E32").Value = '  -7.40%  '

$ws.Range("D33").Value = '12.49'
$ws.Range("E33").Value = '  -8.08%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.117'
$ws.Range("E34").Value = '  -7.06%  '

$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = '65.29'
$ws.Range("E35").Value = '  -4.46%  '

$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").Value = '43.24'
$ws.Range("E36").Value = '  -16.92%  '

$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").Value = '599.80'
$ws.Range("E37").Value = '  -3.55%  '

$ws.Range("D38").Value = '0.0₃0894'
$ws.Range("E38").Value = '  -8.86%  '

$ws.Range("E39").Value = '  +0.05%  '

$ws.Range("D40").Value = '0.398'
$ws.Range("E40").Value = '  -5.16%  '

$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.19%  '

$ws.Range("E42").Value = '  -4.22%  '

$ws.Range("D43").Value = '2.77'
$ws.Range("E43").Value = '  +5.29%  '

$ws.Range("D44").Value = '2.99'
$ws.Range("E44").Value = '  -9.24%  '

$ws.Range("E45").Value = '  -8.04%  '

$ws.Range("D46").Value = '0.0434'
$ws.Range("E46").Value = '  -7.52%  '

$ws.Range("D47").Value = '9.24'
$ws.Range("E47").Value = '  -9.94%  '

$ws.Range("D48").Value = '2.805.66'
$ws.Range("E48").Value = '  -2.10%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '0.133'
$ws.Range("E49").Value = '  -7.10%  '

$ws.Range("B50").Value = 'WEMIXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").Value = '2.69'
$ws.Range("E50").Value = '  -5.09%  '

$ws.Range("D51").Value = '3.11'
$ws.Range("E51").Value = '  -6.74%  '
